$d = $word.ActiveDocument

# The "Range-based Parking Provision Standards" figure is currently embedded
# as an inline picture. Replace it with a hyperlink that points at the
# image's URL on the URA website (the picture -> link swap described by the
# commit).
$url = "https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Industrial/Range_Based_Car_Parking_Standard.jpg?h=100%25&w=100%25"

if ($d.InlineShapes.Count -ge 1) {
    $shape = $d.InlineShapes(1)
    $shapeRange = $shape.Range
    $shape.Delete()
    $d.Hyperlinks.Add($shapeRange, $url, [ref]"", [ref]"", $url) | Out-Null
}
